$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 219 is a blank "separator" row, matching the style of the other
# separator rows already present in the sheet (e.g. row 209).
$ws.Range("A209:C209").Copy($ws.Range("A219:C219"))

# New timesheet entries for 27-02-2020 (rows 220-229), following the
# existing Timestamp / Task / Location column layout.
$data = @(
  @("Feb 27 10:00 to 11:00", "Working with example of creating sample django site for PAM data exploration.", "Infimetrics"),
  @("Feb 27 11:00 to 12:00", "Created feature engineering class and used it to create features", "Infimetrics"),
  @("Feb 27 12:00 to 13:00", "Done modelling on data, used multi output regressor.", "Infimetrics"),
  @("Feb 27 13:00 to 13:30", "Working on tuning model", "Infimetrics"),
  @("Feb 27 13:30 to 14:00", "Lunch", "Infimetrics"),
  @("Feb 27 14:00 to 15:00", "Ml phase done", "Infimetrics"),
  @("Feb 27 15:00 to 16:00", "Working on django backend", "Infimetrics"),
  @("Feb 27 16:00 to 17:00", "Done with sample example, using post making predictions.", "Infimetrics"),
  @("Feb 27 17:00 to 18:00", "Working on deploying classification problem.", "Infimetrics"),
  @("Feb 27 18:00 to 19:00", "Done demo ml example, deployed in django successfully", "Infimetrics")
)

$row = 220
foreach ($entry in $data) {
  $ws.Cells.Item($row, 1).Value = $entry[0]
  $ws.Cells.Item($row, 2).Value = $entry[1]
  $ws.Cells.Item($row, 3).Value = $entry[2]
  $row = $row + 1
}

# Update the saved view state to match where Excel was scrolled/selected
# after the edit (scrolled so row 210 is at the top, active cell B224).
$win = $excel.ActiveWindow
$win.ScrollRow = 210
$win.ScrollColumn = 1
$ws.Range("B224").Select()
